# Add a new "metadata" worksheet after the existing "data" sheet and
# populate it with panel metadata, matching the target OOXML diff.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "metadata"

# Match the page margins used on the "data" sheet (PageSetup margins are in
# points: 1 inch = 72 points).
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Header row (row 1).
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

# Reuse the same bold/bordered header style already used by the "data"
# sheet's header row (B1:F1), so no duplicate style entries are minted.
$dataSheet.Range("B1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats

# Data row (row 2).
$ws.Range("A2").Value = 0
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats - same header style as A2 in the diff

$ws.Range("B2").Value = "Common deletion and duplication syndromes"
$ws.Range("C2").Value = 3443

# data_version looks numeric ("0.137") but must be stored as text, so force
# a text number format before assigning it, then drop back to the default
# (unstyled) cell format so no explicit style index is left on the cell.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0.137"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "2020-12-07T23:28:46.216350Z"
$ws.Range("F2").Value = "2021-10-05 14:33:30.940740"
$ws.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3443/?format=json"

$excel.CutCopyMode = $false

Write-Output "metadata sheet added"
